$d = $word.ActiveDocument

# The document contains three occurrences of "Chad" that must become "ChAD"
# (i.e. the substring "ad" -> "AD"), each time turning the single run that
# carried the word into three runs: "Ch" / "AD" / <rest of original run text>,
# all sharing identical run formatting (Times New Roman).
#
# The host engine auto-merges adjacent runs that end up with identical
# formatting whenever a Range.Text assignment touches a paragraph. Toggling
# a (no-op) formatting property around the text assignment forces the edited
# span to remain a distinct run even after the formatting is reverted,
# matching the run layout produced by a real Word edit.

$count = 0
while ($true) {
  $full = $d.Content
  $fullText = $full.Text
  $idx = $fullText.IndexOf("Chad")
  if ($idx -lt 0) { break }

  # Only the trailing "ad" of "Chad" is rewritten to "AD"; this both
  # changes the text and introduces the run boundary right after "Ch",
  # leaving the remainder of the original run (e.g. " service ...") as a
  # third, separately-split run.
  $sub = $d.Range($idx + 2, $idx + 4)
  $sub.Bold = 1
  $sub.Text = "AD"
  $sub.Bold = 0

  $count = $count + 1
  if ($count -gt 10) { break }
}
